# Auto-generated Excel COM-interop script applying crypto price/volume update
# Commit: Updated cryptos list on Mon Dec  4 15:47:24 UTC 2023 with GitHub Actions
#
# Notes:
#  - Column D holds "numbers" that are really text labels (e.g. thousand-dot
#    formatted prices like "41.446.03", or plain decimals like "229.79"). The
#    sheet stores Price/Volume as text (inline/shared strings), never as real
#    numeric cells, so assignments that "look like a number" to Excel are
#    forced to text with a leading apostrophe and then ClearFormats() strips
#    the resulting quote-prefix cell style back off so no stray formatting is
#    left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.446.03"
$ws.Range("E2").Value = "  +4.26%  "

# Row 3
$ws.Range("D3").Value = "2.217.33"
$ws.Range("E3").Value = "  +2.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'229.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.14%  "

# Row 6
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("D7").Value = "'61.05"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.52%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.400"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.38%  "

# Row 10
$ws.Range("D10").Value = "'57.99"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.70%  "

# Row 11
$ws.Range("D11").Value = "'0.0892"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.57%  "

# Row 12
$ws.Range("E12").Value = "  -0.34%  "

# Row 13
$ws.Range("D13").Value = "2.541.84"
$ws.Range("E13").Value = "  +2.52%  "

# Row 14
$ws.Range("D14").Value = "'15.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.24%  "

# Row 15
$ws.Range("D15").Value = "'21.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.46%  "

# Row 16
$ws.Range("D16").Value = "'0.793"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("D17").Value = "'5.52"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.54%  "

# Row 18
$ws.Range("D18").Value = "2.203.35"
$ws.Range("E18").Value = "  +1.94%  "

# Row 19
$ws.Range("D19").Value = "41.328.53"
$ws.Range("E19").Value = "  +4.20%  "

# Row 20
$ws.Range("D20").Value = "'72.59"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.34%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0894"
$ws.Range("E21").Value = "  +5.79%  "

# Row 22
$ws.Range("D22").Value = "'6.03"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("D23").Value = "'252.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.84%  "

# Row 24
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.43%  "

# Row 27
$ws.Range("E27").Value = "  -0.20%  "

# Row 28
$ws.Range("D28").Value = "'167.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.04%  "

# Row 29
$ws.Range("D29").Value = "'0.140"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.88%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'19.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("D32").Value = "'2.54"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.38%  "

# Row 33
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("E34").Value = "  +6.26%  "

# Row 35
$ws.Range("D35").Value = "'4.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.56%  "

# Row 36
$ws.Range("D36").Value = "'0.0619"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.67"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.13%  "

# Row 38
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'6.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.63%  "

# Row 39
$ws.Range("E39").Value = "  -2.37%  "

# Row 40
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").Value = "'0.000236"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +28.22%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0236"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.02%  "

# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.08%  "

# Row 44
$ws.Range("D44").Value = "'8.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.38%  "

# Row 45
$ws.Range("D45").Value = "'0.0974"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.06%  "

# Row 46
$ws.Range("D46").Value = "'98.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.68%  "

# Row 47
$ws.Range("D47").Value = "'1.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.05%  "

# Row 48
$ws.Range("D48").Value = "1.463.08"
$ws.Range("E48").Value = "  -3.46%  "

# Row 49
$ws.Range("D49").Value = "'16.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.07%  "

# Row 50
$ws.Range("E50").Value = "  -0.79%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.50%  "
